# Applies the "Error Calculations and Plots" edit:
#  - Updates several individual missing-value cells in columns B/D/F
#    (swapping which cells are treated as missing/imputed).
#  - Removes two rows that no longer belong in the cleaned dataset:
#      the row with ID "RM 232" and the row with ID "SC 92".
#    Removing them shifts all subsequent rows up, shrinking the sheet
#    from A1:F35 down to A1:F33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D / F value swaps for rows 2-25 (row numbers unaffected by the later deletes) ---

# Row 2: D2 becomes a real number
$ws.Range("D2").Value = -13.5

# Row 3: D3 becomes missing
$ws.Range("D3").ClearContents()

# Row 4: D4 becomes missing
$ws.Range("D4").ClearContents()

# Row 5: F5 becomes missing
$ws.Range("F5").ClearContents()

# Row 8: F8 becomes a real number
$ws.Range("F8").Value = 17.05

# Row 10: F10 becomes a real number
$ws.Range("F10").Value = 16.43

# Row 11: D11 becomes a real number
$ws.Range("D11").Value = -15.5

# Row 12: F12 becomes missing
$ws.Range("F12").ClearContents()

# Row 13: D13 becomes missing
$ws.Range("D13").ClearContents()

# Row 15: F15 becomes a real number
$ws.Range("F15").Value = 16.2

# Row 18: F18 becomes missing
$ws.Range("F18").ClearContents()

# Row 19: F19 becomes missing
$ws.Range("F19").ClearContents()

# Row 21: D21 becomes a real number
$ws.Range("D21").Value = -14.3

# Row 25: D25 becomes missing, F25 becomes a real number
$ws.Range("D25").ClearContents()
$ws.Range("F25").Value = 16.6

# --- Remove the "RM 232" row entirely (originally row 26) ---
$ws.Rows.Item(26).Delete()

# After the delete above, the row that was "SC 92" (originally row 28)
# has shifted up to row 27. Remove it too, finishing the shift so the
# remaining "SC ..." rows move up by two rows in total.
$ws.Rows.Item(27).Delete()

# --- Additional missing-value swaps among the rows that shifted up ---

# Row 27 (now "SC 101"): F27 becomes a real number
$ws.Range("F27").Value = 17.0

# Row 29 (now "SC 119"): B29 and F29 become missing
$ws.Range("B29").ClearContents()
$ws.Range("F29").ClearContents()

# Row 33 (now "SC 232"): B33 and D33 become real numbers, F33 becomes missing
$ws.Range("B33").Value = -19.5
$ws.Range("D33").Value = -14.1
$ws.Range("F33").ClearContents()
